$d = $word.ActiveDocument
$d.Content.Find.Execute("global_trigger_testcase_inactivity_watchdog", $false, $false, $false, $false, $false, $true, 1, $false, "global_trigger_activity_watchdog", 2)
